$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (638) down
# across the new rows (639-653) for all columns A:N so the new cells get
# the same cell styles (s attribute) as the rest of the data.
$ws.Range("A638:N638").Copy()
$ws.Range("A639:N653").PasteSpecial(-4122)

# Row 638 has no value (and thus no style) in column N, so the paste above
# left column N without a style for the new rows. Pull the column-N style
# from an earlier row (row 2) that does have it, just for that column.
$ws.Range("N2").Copy()
$ws.Range("N639:N653").PasteSpecial(-4122)

# Each response row only has an answer in column M OR column N (single-choice
# question), never both. Remove the formatting-only placeholder cell that was
# produced by the paste above for whichever of M/N has no answer in this batch,
# so it matches the source data exactly (no stray empty cell).
$ws.Range("N639").Clear()
$ws.Range("M640").Clear()
$ws.Range("N641").Clear()
$ws.Range("N642").Clear()
$ws.Range("N643").Clear()
$ws.Range("N644").Clear()
$ws.Range("M645").Clear()
$ws.Range("M646").Clear()
$ws.Range("N647").Clear()
$ws.Range("M648").Clear()
$ws.Range("M649").Clear()
$ws.Range("M650").Clear()
$ws.Range("M651").Clear()
$ws.Range("M652").Clear()
$ws.Range("N653").Clear()

# Populate the new survey response rows (639-653).
$ws.Range("A639").Value = 45192.67608010417
$ws.Range("B639").Value = "mjh8429@naver.com"
$ws.Range("C639").Value = "디지털미디어콘텐츠"
$ws.Range("D639").Value = 20192535
$ws.Range("E639").Value = "민지혜"
$ws.Range("F639").Value = "'78:22"
$ws.Range("G639").Value = 0.15
$ws.Range("H639").Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Range("I639").Value = "779만 명"
$ws.Range("J639").Value = 0.151
$ws.Range("K639").Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Range("L639").Value = "Red"
$ws.Range("M639").Value = "모름/무응답"
$ws.Range("A640").Value = 45192.67764202546
$ws.Range("B640").Value = "041030top@naver.com"
$ws.Range("C640").Value = "데이터사이언스"
$ws.Range("D640").Value = 20233257
$ws.Range("E640").Value = "최영국"
$ws.Range("F640").Value = "'74:26"
$ws.Range("G640").Value = 0.2
$ws.Range("H640").Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Range("I640").Value = "952만 명"
$ws.Range("J640").Value = 0.059
$ws.Range("K640").Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Range("L640").Value = "Black"
$ws.Range("N640").Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A641").Value = 45192.68418751158
$ws.Range("B641").Value = "ekgus0916@naver.com"
$ws.Range("C641").Value = "체육학과"
$ws.Range("D641").Value = 20217125
$ws.Range("E641").Value = "김다현"
$ws.Range("F641").Value = "'74:26"
$ws.Range("G641").Value = 0.1
$ws.Range("H641").Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Range("I641").Value = "779만 명"
$ws.Range("J641").Value = 0.151
$ws.Range("K641").Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Range("L641").Value = "Red"
$ws.Range("M641").Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A642").Value = 45192.71004814815
$ws.Range("B642").Value = "jinyoung05503@gmail.com"
$ws.Range("C642").Value = "생명과학과"
$ws.Range("D642").Value = 20233534
$ws.Range("E642").Value = "이진영"
$ws.Range("F642").Value = "'75:25"
$ws.Range("G642").Value = 0.2
$ws.Range("H642").Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Range("I642").Value = "779만 명"
$ws.Range("J642").Value = 0.151
$ws.Range("K642").Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Range("L642").Value = "Red"
$ws.Range("M642").Value = "모름/무응답"
$ws.Range("A643").Value = 45192.72864728009
$ws.Range("B643").Value = "soccert71@naver.com"
$ws.Range("C643").Value = "심리학과"
$ws.Range("D643").Value = 20182125
$ws.Range("E643").Value = "육정민"
$ws.Range("F643").Value = "'74:26"
$ws.Range("G643").Value = 0.2
$ws.Range("H643").Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Range("I643").Value = "952만 명"
$ws.Range("J643").Value = 0.059
$ws.Range("K643").Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Range("L643").Value = "Red"
$ws.Range("M643").Value = "모름/무응답"
$ws.Range("A644").Value = 45192.72900384259
$ws.Range("B644").Value = "han7434370@naver.com"
$ws.Range("C644").Value = "체육학과"
$ws.Range("D644").Value = 20224152
$ws.Range("E644").Value = "한진우"
$ws.Range("F644").Value = "'77:23"
$ws.Range("G644").Value = 0.2
$ws.Range("H644").Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Range("I644").Value = "166만 명"
$ws.Range("J644").Value = 0.151
$ws.Range("K644").Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Range("L644").Value = "Red"
$ws.Range("M644").Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A645").Value = 45192.73253275463
$ws.Range("B645").Value = "milovany03@gmail.com"
$ws.Range("C645").Value = "사회학과"
$ws.Range("D645").Value = 20202223
$ws.Range("E645").Value = "박진옥"
$ws.Range("F645").Value = "'74:26"
$ws.Range("G645").Value = 0.2
$ws.Range("H645").Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Range("I645").Value = "952만 명"
$ws.Range("J645").Value = 0.059
$ws.Range("K645").Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Range("L645").Value = "Black"
$ws.Range("N645").Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A646").Value = 45192.75874024305
$ws.Range("B646").Value = "catboom5329@gmail.com"
$ws.Range("C646").Value = "체육학과"
$ws.Range("D646").Value = 20234121
$ws.Range("E646").Value = "박주현"
$ws.Range("F646").Value = "'78:22"
$ws.Range("G646").Value = 0.25
$ws.Range("H646").Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Range("I646").Value = "952만 명"
$ws.Range("J646").Value = 0.151
$ws.Range("K646").Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Range("L646").Value = "Black"
$ws.Range("N646").Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A647").Value = 45192.7610846875
$ws.Range("B647").Value = "yeshin05@naver.com"
$ws.Range("C647").Value = "미래융합스쿨"
$ws.Range("D647").Value = 20236639
$ws.Range("E647").Value = "최예원"
$ws.Range("F647").Value = "'76:24"
$ws.Range("G647").Value = 0.25
$ws.Range("H647").Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Range("I647").Value = "779만 명"
$ws.Range("J647").Value = 0.151
$ws.Range("K647").Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Range("L647").Value = "Red"
$ws.Range("M647").Value = "모름/무응답"
$ws.Range("A648").Value = 45192.773353344906
$ws.Range("B648").Value = "ind_b3@naver.com"
$ws.Range("C648").Value = "미디어스쿨"
$ws.Range("D648").Value = 20232523
$ws.Range("E648").Value = "김지안"
$ws.Range("F648").Value = "'76:24"
$ws.Range("G648").Value = 0.2
$ws.Range("H648").Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Range("I648").Value = "779만 명"
$ws.Range("J648").Value = 0.059
$ws.Range("K648").Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Range("L648").Value = "Black"
$ws.Range("N648").Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A649").Value = 45192.77901452546
$ws.Range("B649").Value = "jamesjm0612@gmail.com"
$ws.Range("C649").Value = "영어영문학과"
$ws.Range("D649").Value = 20231231
$ws.Range("E649").Value = "정재민"
$ws.Range("F649").Value = "'75:25"
$ws.Range("G649").Value = 0.15
$ws.Range("H649").Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Range("I649").Value = "166만 명"
$ws.Range("J649").Value = 0.374
$ws.Range("K649").Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Range("L649").Value = "Black"
$ws.Range("N649").Value = "모름/무응답"
$ws.Range("A650").Value = 45192.787771874995
$ws.Range("B650").Value = "kkhe2370@naver.com"
$ws.Range("C650").Value = "광고홍보학과"
$ws.Range("D650").Value = 20202638
$ws.Range("E650").Value = "전혜린"
$ws.Range("F650").Value = "'77:23"
$ws.Range("G650").Value = 0.1
$ws.Range("H650").Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Range("I650").Value = "166만 명"
$ws.Range("J650").Value = 0.151
$ws.Range("K650").Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Range("L650").Value = "Black"
$ws.Range("N650").Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A651").Value = 45192.790401087965
$ws.Range("B651").Value = "jehuncho03@gmail.com"
$ws.Range("C651").Value = "글로벌비즈니스"
$ws.Range("D651").Value = 20226425
$ws.Range("E651").Value = "조제헌"
$ws.Range("F651").Value = "'74:26"
$ws.Range("G651").Value = 0.3
$ws.Range("H651").Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Range("I651").Value = "779만 명"
$ws.Range("J651").Value = 0.059
$ws.Range("K651").Value = "중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다"
$ws.Range("L651").Value = "Black"
$ws.Range("N651").Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A652").Value = 45192.79356601852
$ws.Range("B652").Value = "dms95123@naver.com"
$ws.Range("C652").Value = "사회복지학부"
$ws.Range("D652").Value = 20232317
$ws.Range("E652").Value = "김은별"
$ws.Range("F652").Value = "'74:26"
$ws.Range("G652").Value = 0.2
$ws.Range("H652").Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Range("I652").Value = "952만 명"
$ws.Range("J652").Value = 0.059
$ws.Range("K652").Value = "중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다"
$ws.Range("L652").Value = "Black"
$ws.Range("N652").Value = "모름/무응답"
$ws.Range("A653").Value = 45192.804232407405
$ws.Range("B653").Value = "kangsamy2@gmail.com"
$ws.Range("C653").Value = "사회복지학부"
$ws.Range("D653").Value = 20232302
$ws.Range("E653").Value = "강새미"
$ws.Range("F653").Value = "'76:24"
$ws.Range("G653").Value = 0.2
$ws.Range("H653").Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Range("I653").Value = "779만 명"
$ws.Range("J653").Value = 0.151
$ws.Range("K653").Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Range("L653").Value = "Red"
$ws.Range("M653").Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# Reflect the final selection state recorded in the source workbook.
$ws.Range("F658").Select()

